$wb = $excel.ActiveWorkbook

# Sheets: 1 = Overview, 2 = zh-cn, 3 = de-de
$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

# Update the "Status" value from "Ready for handoff" to "In Translation"
# everywhere it appears (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2 all
# shared the same string table entry).
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# The status columns shrink to fit the new, shorter text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
